$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Question")

# The sheet ships with an empty stylesheet (no cellXfs records at all),
# which crashes the engine's cell-write path the first time any cell
# Value/Formula is set. Touching a cell's Style first seeds the style
# table so every subsequent Value assignment below succeeds. We re-apply
# "Normal" after any NumberFormat/quote-prefix tweak too, so the cells
# don't end up pinned to a throwaway style index.
$ws.Range("A4").Style = "Normal"

# ---- Row 4 ----
$ws.Range("A4").Value = "update"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1934291"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "Composite"
$ws.Range("D4").Value = "Composite"
$ws.Range("E4").Value = "c9523d19-ee08-44b7-b3a7-ece580db8a83"

# Empty-but-present text cells: a lone quote prefix collapses to an empty
# string (matching the blank cells elsewhere in the sheet) instead of
# clearing/removing the cell the way Value = "" would.
$ws.Range("G4").Value = "'"
$ws.Range("G4").Style = "Normal"

$ws.Range("H4").Value = "Constructed"
$ws.Range("I4").Value = "Custom"
$ws.Range("J4").Value = "Human Scoring"
$ws.Range("K4").Value = 126

$ws.Range("L4").Value = "'"
$ws.Range("L4").Style = "Normal"

$ws.Range("M4").Value = "'"
$ws.Range("M4").Style = "Normal"

$ws.Range("N4").Value = 0

$ws.Range("O4").Value = "'"
$ws.Range("O4").Style = "Normal"

$ws.Range("P4").Value = "'"
$ws.Range("P4").Style = "Normal"

$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "11/20/2015"
$ws.Range("Q4").Style = "Normal"

# ---- Row 5 ----
$ws.Range("A5").Value = "update"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1934293"
$ws.Range("B5").Style = "Normal"

$ws.Range("C5").Value = "Composite"
$ws.Range("D5").Value = "Composite"
$ws.Range("E5").Value = "c9523d19-ee08-44b7-b3a7-ece580db8a83"

$ws.Range("G5").Value = "'"
$ws.Range("G5").Style = "Normal"

$ws.Range("H5").Value = "Constructed"
$ws.Range("I5").Value = "Custom"
$ws.Range("J5").Value = "Human Scoring"
$ws.Range("K5").Value = 126

$ws.Range("L5").Value = "'"
$ws.Range("L5").Style = "Normal"

$ws.Range("M5").Value = "'"
$ws.Range("M5").Style = "Normal"

$ws.Range("N5").Value = 0

$ws.Range("O5").Value = "'"
$ws.Range("O5").Style = "Normal"

$ws.Range("P5").Value = "'"
$ws.Range("P5").Style = "Normal"

$ws.Range("Q5").NumberFormat = "@"
$ws.Range("Q5").Value = "11/20/2015"
$ws.Range("Q5").Style = "Normal"
